# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K" - strikeouts) with freshly-computed values,
# replacing the old "Strike#" counts that used to live in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value for column G (rows 2..70 correspond to data rows)
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 2
    6  = 2
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 3
    15 = 2
    16 = 0
    17 = 0
    18 = 2
    19 = 2
    20 = 1
    21 = 0
    22 = 0
    23 = 0
    24 = 2
    25 = 0
    26 = 0
    27 = 1
    28 = 2
    29 = 1
    30 = 2
    31 = 0
    32 = 2
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 1
    38 = 2
    39 = 0
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 0
    45 = 2
    46 = 0
    47 = 1
    48 = 1
    49 = 1
    50 = 1
    51 = 4
    52 = 1
    53 = 0
    54 = 0
    55 = 2
    56 = 0
    57 = 0
    58 = 2
    59 = 1
    60 = 0
    61 = 2
    62 = 2
    63 = 0
    64 = 2
    65 = 1
    66 = 0
    67 = 1
    68 = 2
    69 = 2
    70 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
